$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$value = "01/06/2016 (Olivier)"

$cells = @("C5","D5","C7","C10","C11","C12","F14","C16","C27","C30","C31","C34","C38","C40","C42","D42","E42","F42","C44","C53","C54")
foreach ($addr in $cells) {
    $ws.Range($addr).Value = $value
}

# Rows 14 and 42 are now fully complete (Total = 5) -> apply the same
# "complete row" highlight already used on A6 / A17 / A21 to column A.
$ws.Range("A6").Copy()
$ws.Range("A14").PasteSpecial(-4122)
$ws.Range("A42").PasteSpecial(-4122)

# C31 gets the same "filled" highlight already used on cells like C6.
$ws.Range("C6").Copy()
$ws.Range("C31").PasteSpecial(-4122)

# Selection moved to C3 and top-left cell reset (per sheetView change in the diff)
$ws.Range("C3").Select()
